# The underlying source rows (title/timestamp/uri triples) got reshuffled
# into a new row order - add one json for time bucket analysis caused the
# 8 existing entries to be re-emitted in a different sequence. Column
# headers (row 1) and the "historical distance"/"time bucket" columns
# (C, D - always "unknown") are unaffected; only A (title), B (timestamp)
# and E (uri) text per row changes. Existing hyperlinks (E2:E9, already
# wired to their own relationship ids) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Presidential Ratings"
$ws.Range("B2").Value = "1-01-01T00:00:00UTC"
$ws.Range("E2").Value = "https://insideelections.com/ratings/president"

$ws.Range("A3").Value = "Find Your Local League"
$ws.Range("B3").Value = "1-01-01T00:00:00UTC"
$ws.Range("E3").Value = "https://www.lwv.org/local-leagues/find-local-league"

$ws.Range("A4").Value = "Polls 2020-11-02 (smaller states)"
$ws.Range("B4").Value = "2020-11-02T00:00:00UTC"
$ws.Range("E4").Value = "https://www.swayable.com/polls/2020-11-02-small.html"

$ws.Range("A5").Value = "Maryland: Election Tools, Deadlines, Dates, Rules, and Links"
$ws.Range("B5").Value = "1-01-01T00:00:00UTC"
$ws.Range("E5").Value = "https://www.vote.org/state/maryland/"

$ws.Range("A6").Value = "Biden dominates the electoral map, but here's how the race could tighten"
$ws.Range("B6").Value = "2020-08-06T13:13:00UTC"
$ws.Range("E6").Value = "https://www.nbcnews.com/politics/meet-the-press/biden-dominates-electoral-map-here-s-how-race-could-tighten-n1236001"

$ws.Range("A7").Value = "Maryland postpones primary, shifts special election to mail voting over coronavirus"
$ws.Range("B7").Value = "2020-03-17T11:11:00UTC"
$ws.Range("E7").Value = "https://www.politico.com/news/2020/03/17/maryland-postpones-april-28-primary-election-over-coronavirus-133776"

$ws.Range("A8").Value = "2020 Election Forecast"
$ws.Range("B8").Value = "2020-08-12T06:30:00UTC"
$ws.Range("E8").Value = "https://projects.fivethirtyeight.com/2020-election-forecast/"

$ws.Range("A9").Value = "2020 Electoral Interactive Map"
$ws.Range("B9").Value = "1-01-01T00:00:00UTC"
$ws.Range("E9").Value = "https://abcnews.go.com/Politics/2020-Electoral-Interactive-Map?basemap=71662160&promoref=brandpromo"
